# Add data for 2022-06-24: update sheet name, label string, and June/Total rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet/tab to reflect the new "through" date.
$ws.Name = "Through 2022-06-16"

# Update the June row label.
$ws.Range("A7").Value = "June (through 06-16)"

# Update June row (row 7) values.
$ws.Range("B7").Value = 8
$ws.Range("C7").Value = 20
$ws.Range("D7").Value = 31
$ws.Range("E7").Value = 32
$ws.Range("F7").Value = 26
$ws.Range("G7").Value = 60
$ws.Range("H7").Value = 54
$ws.Range("I7").Value = 74

# Update Total row (row 8) values.
$ws.Range("B8").Value = 116
$ws.Range("C8").Value = 229
$ws.Range("D8").Value = 347
$ws.Range("E8").Value = 327
$ws.Range("F8").Value = 230
$ws.Range("G8").Value = 418
$ws.Range("H8").Value = 685
$ws.Range("I8").Value = 737
